$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Update the departure date (E2) and arrival date (K2)
$ws.Range("E2").Value = "30/11/2017"
$ws.Range("K2").Value = "01/12/2017"

# Update the arrival hour (L2) to a text value "19h" instead of the raw number 18
$ws.Range("L2").Value = "19h"

# Move the active selection to D5
$ws.Range("D5").Select()
